$wb = $excel.ActiveWorkbook
$czech = $wb.Worksheets.Item("Czech")
$czech.Copy($null, $czech)
$swiss = $wb.Worksheets.Item($wb.Worksheets.Count)
$swiss.Name = "Swiss"
$swiss.Range("B2").Value = "Switzerland Market"
$swiss.Range("B4").Value = "NGC-3476/T2656"

$germany = $wb.Worksheets.Item("Germany")
$germany.Activate()
[void]$germany.Cells.Select()

$swiss.Activate()
[void]$swiss.Range("A10").Select()
